# Add Tan's schedule: rename Sheet3 -> "Tân", populate it as a copy of
# "C.Thắng"'s schedule (values + styles), then re-color a handful of
# cells to reflect Tan's own weekly timetable, and finally make the new
# sheet the active one (matches workbookView activeTab + tabSelected).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# 1. Rename the blank third sheet.
$ws3.Name = "Tân"

# 2. Clone C.Thắng's whole grid (values, shared-string cells and styles)
#    onto the new sheet.
$ws2.Cells.Copy($ws3.Range("A1"))

# 3. Re-paint Tan's own busy/free slots.
#    3a. A few slots that are free for Tan but booked (yellow) for
#        C.Thắng become plain "white" cells (new fill, theme color 2 ->
#        theme="0").
$ws3.Range("C2:C7").Interior.ColorIndex = 6
$ws3.Range("C2:C7").Interior.ThemeColor = 2
$ws3.Range("C16:C21").Interior.ColorIndex = 6
$ws3.Range("C16:C21").Interior.ThemeColor = 2
$ws3.Range("F15").Interior.ColorIndex = 6
$ws3.Range("F15").Interior.ThemeColor = 2
$ws3.Range("F20:F25").Interior.ColorIndex = 6
$ws3.Range("F20:F25").Interior.ThemeColor = 2
$ws3.Range("G15:G25").Interior.ColorIndex = 6
$ws3.Range("G15:G25").Interior.ThemeColor = 2

#    3b. A block of slots that become Tan's own "blue" busy marker (new
#        fill, theme color 5 -> theme="4").
$ws3.Range("B16:B19").Interior.ColorIndex = 6
$ws3.Range("B16:B19").Interior.ThemeColor = 5
$ws3.Range("D16:D19").Interior.ColorIndex = 6
$ws3.Range("D16:D19").Interior.ThemeColor = 5
$ws3.Range("F16:F19").Interior.ColorIndex = 6
$ws3.Range("F16:F19").Interior.ThemeColor = 5

#    3c. A few more slots that become "booked" (existing yellow fill).
$ws3.Range("D10:D12").Interior.ColorIndex = 6
$ws3.Range("E13:E14").Interior.ColorIndex = 6

# 4. Tan's sheet keeps going down to row 25 (C.Thắng's stopped at 19) -
#    extend the grid with blank, white-filled cells.
$ws3.Range("C20:C21").Interior.ColorIndex = 6
$ws3.Range("C20:C21").Interior.ThemeColor = 2

# 5. Restore the selection on C.Thắng (it is no longer the active tab)
#    and set Tan's own selection, then activate Tan's sheet last so it
#    becomes the workbook's active tab/sheet.
$ws2.Range("A1:XFD1048576").Select()
$ws3.Range("F8").Select()
$ws3.Activate()
